$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 20..81 down to 21..82 (working bottom-up so we never overwrite
# a row before it has been read), freeing up row 20 for the new record.
for ($r = 81; $r -ge 20; $r--) {
    for ($col = 1; $col -le 5; $col++) {
        $srcVal = $ws.Cells.Item($r, $col).Value2
        $ws.Cells.Item($r + 1, $col).Value = $srcVal
    }
}

# Fill in the new row 20: Dia=25, total_venda=10636.32, Mes=7, Ano=2025, Periodo=07/2025
$ws.Range("A20").Value = 25
$ws.Range("B20").Value = 10636.32
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 2025
$ws.Range("E20").Value = "07/2025"

# Update the revised total_venda figures for July (rows 2-19 keep their
# original position, only the values changed).
$ws.Range("B2").Value = 18093.98
$ws.Range("B15").Value = 4374.5
$ws.Range("B17").Value = 9220.860000000001
$ws.Range("B18").Value = 12454.83
$ws.Range("B19").Value = 4369

# Update the revised total_venda for June day=12 (now sitting at row 29
# after the insertion above).
$ws.Range("B29").Value = 18059.33
